# Update countries & provincias Spain
# Applies the 30-May-2020 10:10 data refresh to the "Pais" sheet:
#  - timestamp footnote update
#  - numeric refreshes for several existing countries
#  - three country blocks get a new/reinserted row (shifting a couple of
#    rows down) and two adjacent-row swaps, which together account for
#    the shared-string reordering seen in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp footnote -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 10:10"

# --- straightforward numeric refreshes (country stays on same row) ----

# Row 6: Rusia
$ws.Cells.Item(6, 2).Value = 396575
$ws.Cells.Item(6, 3).Value = 8952
$ws.Cells.Item(6, 4).Value = 167469
$ws.Cells.Item(6, 5).Value = 224551
$ws.Cells.Item(6, 7).Value = 181
$ws.Cells.Item(6, 8).Value = 4555

# Row 11: Alemania
$ws.Cells.Item(11, 4).Value = 164900
$ws.Cells.Item(11, 5).Value = 9525

# Row 12: India
$ws.Cells.Item(12, 2).Value = 174020
$ws.Cells.Item(12, 3).Value = 529
$ws.Cells.Item(12, 4).Value = 82676
$ws.Cells.Item(12, 5).Value = 86363
$ws.Cells.Item(12, 7).Value = 1
$ws.Cells.Item(12, 8).Value = 4981

# Row 29: Singapur
$ws.Cells.Item(29, 2).Value = 34366
$ws.Cells.Item(29, 3).Value = 506
$ws.Cells.Item(29, 5).Value = 14712

# Row 39: Polonia
$ws.Cells.Item(39, 4).Value = 11016
$ws.Cells.Item(39, 5).Value = 11088

# Row 57: Chequia
$ws.Cells.Item(57, 2).Value = 9200
$ws.Cells.Item(57, 3).Value = 4
$ws.Cells.Item(57, 4).Value = 6502
$ws.Cells.Item(57, 5).Value = 2379

# Row 61: Noruega
$ws.Cells.Item(61, 2).Value = 8425
$ws.Cells.Item(61, 3).Value = 3
$ws.Cells.Item(61, 5).Value = 462

# Row 62: Moldavia
$ws.Cells.Item(62, 4).Value = 4455
$ws.Cells.Item(62, 5).Value = 3152
$ws.Cells.Item(62, 7).Value = 1
$ws.Cells.Item(62, 8).Value = 289

# Row 88: El Salvador
$ws.Cells.Item(88, 2).Value = 2395
$ws.Cells.Item(88, 3).Value = 117
$ws.Cells.Item(88, 4).Value = 1026
$ws.Cells.Item(88, 5).Value = 1325
$ws.Cells.Item(88, 7).Value = 2
$ws.Cells.Item(88, 8).Value = 44

# Row 92: Estonia
$ws.Cells.Item(92, 2).Value = 1865
$ws.Cells.Item(92, 3).Value = 6
$ws.Cells.Item(92, 4).Value = 1622
$ws.Cells.Item(92, 5).Value = 176

# Row 102: Eslovaquia
$ws.Cells.Item(102, 2).Value = 1521
$ws.Cells.Item(102, 3).Value = 1
$ws.Cells.Item(102, 4).Value = 1356
$ws.Cells.Item(102, 5).Value = 137

# --- "Consejo Danes para los Refugiados" moves ahead of "Republica de
#     Yibuti"/"Grecia" with refreshed numbers; Yibuti & Grecia's rows
#     shift down one row each, carrying their prior values along. -----

$ws.Cells.Item(81, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(81, 2).Value = 2966
$ws.Cells.Item(81, 3).Value = 133
$ws.Cells.Item(81, 4).Value = 428
$ws.Cells.Item(81, 5).Value = 2469
$ws.Cells.Item(81, 8).Value = 69

$ws.Cells.Item(82, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(82, 2).Value = 2914
$ws.Cells.Item(82, 4).Value = 1241
$ws.Cells.Item(82, 5).Value = 1653
$ws.Cells.Item(82, 8).Value = 20

$ws.Cells.Item(83, 1).Value = "Grecia"
$ws.Cells.Item(83, 2).Value = 2909
$ws.Cells.Item(83, 4).Value = 1374
$ws.Cells.Item(83, 5).Value = 1360
$ws.Cells.Item(83, 8).Value = 175

# --- adjacent-row swaps (matching stats already equal in B/C/E/F/G) ---

# Rows 198/199: Fiyi <-> Curazao
$ws.Cells.Item(198, 1).Value = "Curazao"
$ws.Cells.Item(198, 4).Value = 14
$ws.Cells.Item(198, 8).Value = 1

$ws.Cells.Item(199, 1).Value = "Fiyi"
$ws.Cells.Item(199, 4).Value = 15
$ws.Cells.Item(199, 8).Value = 0

# Rows 210/211: Seychelles <-> Montserrat
$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1

$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0

# Rows 213/214: Papua Nueva Guinea <-> Islas Virgenes Britanicas
$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 8).Value = 1

$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 8).Value = 0
